$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-07-13"

$ws.Range("A8").Value = "July (through 07-13)"

$ws.Range("B8").Value = 17
$ws.Range("C8").Value = 26
$ws.Range("D8").Value = 23
$ws.Range("E8").Value = 30
$ws.Range("F8").Value = 24
$ws.Range("G8").Value = 45
$ws.Range("H8").Value = 63
$ws.Range("I8").Value = 71

$ws.Range("B9").Value = 142
$ws.Range("C9").Value = 274
$ws.Range("D9").Value = 413
$ws.Range("E9").Value = 383
$ws.Range("F9").Value = 275
$ws.Range("G9").Value = 517
$ws.Range("H9").Value = 823
$ws.Range("I9").Value = 877
